$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Cells.Item(2, 4)
$c.NumberFormat = '@'
$c.Value = '59.528.21'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E2').Value = '  +0.44%  '
$c = $ws.Cells.Item(3, 4)
$c.NumberFormat = '@'
$c.Value = '2.643.48'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E3').Value = '  +1.54%  '
$ws.Range('E4').Value = '  -0.05%  '
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = '@'
$c.Value = '537.19'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  -0.05%  '
$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = '@'
$c.Value = '145.32'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  +3.57%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('E8').Value = '  +0.74%  '
$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = '@'
$c.Value = '6.68'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E9').Value = '  +3.59%  '
$ws.Range('E10').Value = '  -0.20%  '
$ws.Range('E11').Value = '  +0.82%  '
$ws.Range('E12').Value = '  -0.32%  '
$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = '@'
$c.Value = '3.107.90'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = '@'
$c.Value = '59.451.35'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E14').Value = '  +0.42%  '
$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = '@'
$c.Value = '21.20'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E15').Value = '  +3.30%  '
$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = '@'
$c.Value = '2.665.69'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E16').Value = '  +2.40%  '
$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = '@'
$c.Value = '0.0000135'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E17').Value = '  +0.78%  '
$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = '@'
$c.Value = '339.31'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E18').Value = '  -0.81%  '
$ws.Range('E19').Value = '  +1.29%  '
$ws.Range('E20').Value = '  +3.13%  '
$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = '@'
$c.Value = '6.29'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  -1.46%  '
$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  +0.01%  '
$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = '@'
$c.Value = '66.98'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  -0.68%  '
$ws.Range('E24').Value = '  +2.03%  '
$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = '@'
$c.Value = '0.164'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  -1.09%  '
$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = '@'
$c.Value = '0.998'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  -0.09%  '
$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = '@'
$c.Value = '7.27'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  +1.14%  '
$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = '@'
$c.Value = '0.0₃0744'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E28').Value = '  +1.03%  '
$ws.Range('E30').Value = '  +0.17%  '
$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = '@'
$c.Value = '5.85'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E31').Value = '  +0.41%  '
$ws.Range('E32').Value = '  +0.71%  '
$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = '@'
$c.Value = '151.55'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  +1.60%  '
$ws.Range('E34').Value = '  +1.08%  '
$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = '@'
$c.Value = '1.13'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E35').Value = '  +1.74%  '
$ws.Range('E36').Value = '  +2.68%  '
$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = '@'
$c.Value = '0.834'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E37').Value = '  -0.15%  '
$ws.Range('E38').Value = '  -1.66%  '
$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = '@'
$c.Value = '288.06'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  +5.16%  '
$ws.Range('E40').Value = '  +1.60%  '
$ws.Range('E41').Value = '  -0.04%  '
$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = '@'
$c.Value = '0.606'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  +1.72%  '
$ws.Range('E43').Value = '  +0.08%  '
$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = '@'
$c.Value = '19.31'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  +4.07%  '
$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = '@'
$c.Value = '0.0538'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  +2.83%  '
$ws.Range('E46').Value = '  -1.55%  '
$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = '@'
$c.Value = '1.969.62'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  +1.38%  '
$ws.Range('E48').Value = '  +1.58%  '
$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = '@'
$c.Value = '4.56'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E49').Value = '  +1.42%  '
$ws.Range('E50').Value = '  +0.28%  '
$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = '@'
$c.Value = '110.77'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  +0.01%  '
